$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 22, pushing the existing row 22 (and everything
# below it) down by one row (row 119 -> 120, dimension A1:R119 -> A1:R120).
$ws.Rows.Item(22).Insert()

# Populate the newly inserted row 22 with a new weekly price record, using the
# same boilerplate fields (market/region/category/etc.) as the other rows in
# this sheet, but with fresh date/volume/price values.
$ws.Cells.Item(22, 1).Value = 8
$ws.Cells.Item(22, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(22, 3).Value = "Coquimbo"
$ws.Cells.Item(22, 4).Value = 44600
$ws.Cells.Item(22, 5).Value = 4
$ws.Cells.Item(22, 6).Value = 100112040
$ws.Cells.Item(22, 7).Value = "Cilantro"
$ws.Cells.Item(22, 8).Value = "Sin especificar"
$ws.Cells.Item(22, 9).Value = "Primera"
$ws.Cells.Item(22, 10).Value = 2400
$ws.Cells.Item(22, 11).Value = 2300
$ws.Cells.Item(22, 12).Value = 2500
$ws.Cells.Item(22, 13).Value = 2400
$ws.Cells.Item(22, 14).Value = "$/atado 1 a 1,5 kilos"
$ws.Cells.Item(22, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(22, 16).Value = 1600
$ws.Cells.Item(22, 17).Value = 1.5
$ws.Cells.Item(22, 18).Value = "Hortaliza"
